$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated dSF (column F) values per row, reflecting repulled data / mean calculation
$updates = @{
    2  = -3
    3  = -7
    4  = 1
    6  = -3
    8  = 3
    9  = -1
    10 = 1
    12 = 3
    13 = 1
    14 = 1
    15 = 1
    16 = -1
    17 = -1
    18 = -2
    19 = 2
    20 = 5
    21 = -2
    22 = -2
    23 = 6
    24 = 2
    25 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
